$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels to use spaced-out, more readable names.
$ws.Range("A1").Value = "Task ID"
$ws.Range("D1").Value = "Execution Time"
$ws.Range("F1").Value = "Resource Requirements"

# Move the active selection to F1 (matches the saved sheetView selection).
$ws.Range("F1").Select()
